# Weekly crime data update for CompStat report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = -66.666666666666
$ws.Range("C15").NumberFormat = 'General'
$ws.Range("C15").Value = "0"
$ws.Range("D15").NumberFormat = 'General'
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = 'General'
$ws.Range("E15").Value = "***.*"
$ws.Range("M15").Value = -57.142857142857
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -23.728813559322
$ws.Range("L16").Value = -21.052631578947
$ws.Range("M16").Value = -52.127659574468
$ws.Range("N16").Value = -86.880466472303
$ws.Range("C17").NumberFormat = 'General'
$ws.Range("C17").Value = "0"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -29.411764705882
$ws.Range("J17").Value = 98
$ws.Range("K17").Value = -20.408163265306
$ws.Range("L17").Value = -25.714285714285
$ws.Range("M17").Value = -3.703703703703
$ws.Range("N17").Value = -76.506024096385
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -26.666666666666
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = -8.333333333333
$ws.Range("N18").Value = -90.406976744186
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 120
$ws.Range("F19").Value = 44
$ws.Range("H19").Value = 22.222222222222
$ws.Range("I19").Value = 182
$ws.Range("J19").Value = 177
$ws.Range("K19").Value = 2.824858757062
$ws.Range("L19").Value = 26.388888888888
$ws.Range("M19").Value = 163.768115942029
$ws.Range("N19").Value = 13.75
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 36
$ws.Range("K20").Value = 12.5
$ws.Range("L20").Value = -29.411764705882
$ws.Range("M20").Value = 44
$ws.Range("N20").Value = -76
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 5.882352941176
$ws.Range("F21").Value = 82
$ws.Range("H21").Value = 2.5
$ws.Range("I21").Value = 381
$ws.Range("J21").Value = 420
$ws.Range("K21").Value = -9.285714285714
$ws.Range("L21").Value = -2.307692307692
$ws.Range("M21").Value = 18.691588785046
$ws.Range("N21").Value = -72.490974729241
$ws.Range("F22").NumberFormat = 'General'
$ws.Range("F22").Value = "0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 11.111111111111
$ws.Range("M22").Value = 66.666666666666
$ws.Range("D23").NumberFormat = 'General'
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = 'General'
$ws.Range("E23").Value = "***.*"
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -10.526315789473
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = -1.666666666666
$ws.Range("I24").Value = 368
$ws.Range("J24").Value = 342
$ws.Range("K24").Value = 7.602339181286
$ws.Range("L24").Value = 4.545454545454
$ws.Range("M24").Value = 148.648648648649
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -57.692307692307
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = -9.302325581395
$ws.Range("L25").Value = -12.359550561797
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 164
$ws.Range("J26").Value = 159
$ws.Range("K26").Value = 3.144654088050
$ws.Range("L26").Value = 15.492957746478
$ws.Range("M26").Value = -29.004329004329
$ws.Range("C27").NumberFormat = 'General'
$ws.Range("C27").Value = "0"
$ws.Range("D27").NumberFormat = 'General'
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = 'General'
$ws.Range("E27").Value = "***.*"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("I28").Value = 16
$ws.Range("K28").Value = 23.076923076923
$ws.Range("L28").Value = 6.666666666666
